$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header A1 from "Time (day)" to "t_stamp"
$ws.Range("A1").Value = "t_stamp"

# Update selection to A2 (matches sheetView selection change in diff)
$ws.Range("A2").Select()
